$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 127
$ws.Range("I28").Value = 154.5
$ws.Range("J28").Value = 17
$ws.Range("K28").Value = 154.5
$ws.Range("L28").Value = 17
$ws.Range("M28").Value = 330.5
$ws.Range("N28").Value = -987
$ws.Range("H54").Value = 7399.6665
$ws.Range("I54").Value = 7399.6665
$ws.Range("K54").Value = 7399.6665
$ws.Range("M54").Value = -6913.6665
$ws.Range("H115").Value = 2324.2
$ws.Range("I115").Value = 2543.6667
$ws.Range("J115").Value = 349
$ws.Range("K115").Value = 7631.000100000001
$ws.Range("L115").Value = 1047
$ws.Range("M115").Value = -6064.000100000001
$ws.Range("N115").Value = -4181
$ws.Range("H116").Value = 3797.3618
$ws.Range("I116").Value = 3483.9375
$ws.Range("J116").Value = 4466
$ws.Range("K116").Value = 3483.9375
$ws.Range("L116").Value = 4466
$ws.Range("M116").Value = -41.9375
$ws.Range("N116").Value = -11350
$ws.Range("H127").Value = 4263.143
$ws.Range("I127").Value = 4607
$ws.Range("J127").Value = 2200
$ws.Range("K127").Value = 13821
$ws.Range("L127").Value = 6600
$ws.Range("M127").Value = -8861
$ws.Range("N127").Value = -16520
$ws.Range("H132").Value = 24391662
$ws.Range("I132").Value = 24391662
$ws.Range("K132").Value = 73174986
$ws.Range("M132").Value = -73172456
$ws.Range("H137").Value = 3197.5625
$ws.Range("J137").Value = 5095.6665
$ws.Range("L137").Value = 15286.9995
$ws.Range("N137").Value = -20386.9995
$ws.Range("H138").Value = 4829.35
$ws.Range("J138").Value = 6691.475
$ws.Range("L138").Value = 20074.425
$ws.Range("N138").Value = -30354.425

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1828.44
$ws.Range("I2").Value = 1843.4348
$ws.Range("K2").Value = 1843.4348
$ws.Range("M2").Value = -1730.4348
$ws.Range("H32").Value = 3058.9805
$ws.Range("I32").Value = 2333.8958
$ws.Range("K32").Value = 2333.8958
$ws.Range("M32").Value = -2046.8958
$ws.Range("H61").Value = 4272.75
$ws.Range("I61").Value = 4180.857
$ws.Range("J61").Value = 4401.4
$ws.Range("K61").Value = 4180.857
$ws.Range("L61").Value = 4401.4
$ws.Range("M61").Value = -3968.857
$ws.Range("N61").Value = -4825.4
$ws.Range("H74").Value = 2415.611
$ws.Range("I74").Value = 2063.7097
$ws.Range("K74").Value = 2063.7097
$ws.Range("M74").Value = -1189.7097
$ws.Range("H77").Value = 2415.611
$ws.Range("I77").Value = 2063.7097
$ws.Range("K77").Value = 10318.5485
$ws.Range("M77").Value = -5950.548499999999
$ws.Range("H109").Value = 50185
$ws.Range("J109").Value = 50185
$ws.Range("L109").Value = 50185
$ws.Range("N109").Value = -52959
$ws.Range("H112").Value = 29999
$ws.Range("J112").Value = 29999
$ws.Range("L112").Value = 29999
$ws.Range("N112").Value = -32953
$ws.Range("H116").Value = 1828.44
$ws.Range("I116").Value = 1843.4348
$ws.Range("K116").Value = 1843.4348
$ws.Range("M116").Value = 450.5652
$ws.Range("H122").Value = 6455891.5
$ws.Range("I122").Value = 6900687
$ws.Range("K122").Value = 20702061
$ws.Range("M122").Value = -20699611
$ws.Range("H132").Value = 6590.958
$ws.Range("I132").Value = 4785.012
$ws.Range("J132").Value = 19082.084
$ws.Range("K132").Value = 14355.036
$ws.Range("L132").Value = 57246.25199999999
$ws.Range("M132").Value = -11825.036
$ws.Range("N132").Value = -62306.25199999999
$ws.Range("H136").Value = 4272.75
$ws.Range("I136").Value = 4180.857
$ws.Range("J136").Value = 4401.4
$ws.Range("K136").Value = 12542.571
$ws.Range("L136").Value = 13204.2
$ws.Range("M136").Value = -9992.571
$ws.Range("N136").Value = -18304.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1828.44
$ws.Range("I3").Value = 1843.4348
$ws.Range("K3").Value = 1843.4348
$ws.Range("M3").Value = -1729.4348
$ws.Range("H105").Value = 1679.5
$ws.Range("I105").Value = 1712.7778
$ws.Range("J105").Value = 1380
$ws.Range("K105").Value = 1712.7778
$ws.Range("L105").Value = 1380
$ws.Range("M105").Value = 34.22219999999993
$ws.Range("N105").Value = -4874
$ws.Range("H110").Value = 150000
$ws.Range("J110").Value = 150000
$ws.Range("L110").Value = 150000
$ws.Range("N110").Value = -158180
$ws.Range("H132").Value = 51706
$ws.Range("J132").Value = 51706
$ws.Range("L132").Value = 51706
$ws.Range("N132").Value = -61826
$ws.Range("H134").Value = 19232216
$ws.Range("I134").Value = 22728820
$ws.Range("J134").Value = 893.5
$ws.Range("K134").Value = 68186460
$ws.Range("L134").Value = 2680.5
$ws.Range("M134").Value = -68183925
$ws.Range("N134").Value = -7750.5
$ws.Range("H135").Value = 74672.69500000001
$ws.Range("J135").Value = 74672.69500000001
$ws.Range("L135").Value = 74672.69500000001
$ws.Range("N135").Value = -84812.69500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 96.55556
$ws.Range("I4").Value = 96.55556
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 96.55556
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 15.44444
$ws.Range("N4").Value = ""
$ws.Range("H23").Value = 45803.4
$ws.Range("I23").Value = 49499.5
$ws.Range("K23").Value = 49499.5
$ws.Range("M23").Value = -49259.5
$ws.Range("H27").Value = 45803.4
$ws.Range("I27").Value = 49499.5
$ws.Range("K27").Value = 49499.5
$ws.Range("M27").Value = -49307.5
$ws.Range("H31").Value = 2170.5715
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 2170.5715
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 2170.5715
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = -2760.5715
$ws.Range("H34").Value = 2170.5715
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 2170.5715
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 2170.5715
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = -2574.5715
$ws.Range("H58").Value = 3005.3635
$ws.Range("I58").Value = 3337.52
$ws.Range("J58").Value = 1967.375
$ws.Range("K58").Value = 3337.52
$ws.Range("L58").Value = 1967.375
$ws.Range("M58").Value = -3134.52
$ws.Range("N58").Value = -2373.375
$ws.Range("H99").Value = 4284.857
$ws.Range("I99").Value = 3332
$ws.Range("J99").Value = 4999.5
$ws.Range("K99").Value = 3332
$ws.Range("L99").Value = 4999.5
$ws.Range("M99").Value = -1834
$ws.Range("N99").Value = -7995.5
$ws.Range("H126").Value = 4284.857
$ws.Range("I126").Value = 3332
$ws.Range("J126").Value = 4999.5
$ws.Range("K126").Value = 9996
$ws.Range("L126").Value = 14998.5
$ws.Range("M126").Value = -7526
$ws.Range("N126").Value = -19938.5
$ws.Range("H134").Value = 3441.25
$ws.Range("I134").Value = 3441.25
$ws.Range("K134").Value = 10323.75
$ws.Range("M134").Value = -7788.75
$ws.Range("H136").Value = 3005.3635
$ws.Range("I136").Value = 3337.52
$ws.Range("J136").Value = 1967.375
$ws.Range("K136").Value = 10012.56
$ws.Range("L136").Value = 5902.125
$ws.Range("M136").Value = -7462.559999999999
$ws.Range("N136").Value = -11002.125
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""
$ws.Range("H140").Value = 79333.336
$ws.Range("J140").Value = 79333.336
$ws.Range("L140").Value = 79333.336
$ws.Range("N140").Value = -89693.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 10474.75
$ws.Range("I68").Value = 10474.75
$ws.Range("K68").Value = 31424.25
$ws.Range("M68").Value = -30613.25
$ws.Range("H71").Value = 10474.75
$ws.Range("I71").Value = 10474.75
$ws.Range("K71").Value = 94272.75
$ws.Range("M71").Value = -90216.75
$ws.Range("H107").Value = 1123.3334
$ws.Range("J107").Value = 300
$ws.Range("L107").Value = 900
$ws.Range("N107").Value = -4740
$ws.Range("H132").Value = 1060.5555
$ws.Range("I132").Value = 923.8823
$ws.Range("J132").Value = 1292.9
$ws.Range("K132").Value = 8314.940699999999
$ws.Range("L132").Value = 11636.1
$ws.Range("M132").Value = -5784.940699999999
$ws.Range("N132").Value = -16696.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3181.3794
$ws.Range("I132").Value = 2089.4285
$ws.Range("J132").Value = 4200.533
$ws.Range("K132").Value = 6268.2855
$ws.Range("L132").Value = 12601.599
$ws.Range("M132").Value = -3738.2855
$ws.Range("N132").Value = -17661.599

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
$ws.Range("H122").Value = 5142.923
$ws.Range("J122").Value = 6388.222
$ws.Range("L122").Value = 19164.666
$ws.Range("N122").Value = -24064.666
$ws.Range("H140").Value = 79818.27
$ws.Range("J140").Value = 79818.27
$ws.Range("L140").Value = 79818.27
$ws.Range("N140").Value = -90178.27

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 69992
$ws.Range("J46").Value = 69992
$ws.Range("L46").Value = 69992
$ws.Range("N46").Value = -70454
$ws.Range("H132").Value = 1424.9706
$ws.Range("I132").Value = 1387.9642
$ws.Range("K132").Value = 4163.892599999999
$ws.Range("M132").Value = -1633.892599999999
$ws.Range("H134").Value = 69992
$ws.Range("J134").Value = 69992
$ws.Range("L134").Value = 209976
$ws.Range("N134").Value = -215046
